# Opencart_TestCases.xlsx -- record actual run results for the
# still-blank "Result" cells, reusing the formatting that is already used
# for the same status elsewhere in the workbook (so the cell picks up the
# existing "Passed" / "Not run" look instead of a brand new style).

$wb = $excel.ActiveWorkbook

$login    = $wb.Worksheets.Item("Login")
$register = $wb.Worksheets.Item("Register")
$pdp      = $wb.Worksheets.Item("Product Display Page")

# Existing, already-styled cells we can clone the format from.
$passedTemplate = $login.Range("G22")   # fillId 4 ("Passed") look
$notRunTemplate = $pdp.Range("H16")     # fillId 6 ("Not run") look

function Set-ResultCell($range, $text, $template) {
    $template.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null   # xlPasteFormats -- formatting only
    $range.Value = $text
}

# --- Login sheet: fill in the two outstanding results + the last "Not run"
Set-ResultCell $login.Range("G23") "Passed"  $passedTemplate
Set-ResultCell $login.Range("G24") "Passed"  $passedTemplate
Set-ResultCell $login.Range("G25") "Not run" $notRunTemplate

# --- Register sheet: every still-blank result becomes "Not run"
foreach ($r in 14, 16, 19, 22, 23, 24, 25) {
    Set-ResultCell $register.Range("G$r") "Not run" $notRunTemplate
}

$excel.CutCopyMode = 0

# --- Window/selection bookkeeping, restoring Login as the active tab last
$register.Activate()
$register.Range("G14").Select() | Out-Null

$login.Activate()
$login.Application.ActiveWindow.Zoom = 80
$login.Range("K24").Select() | Out-Null
